# Add team record columns (Wins, Losses, Ties) to the KCR_2013 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new column headers
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the formatting from the existing header cell (AC1) onto the new
# header cells so they pick up the same bold/centered/bordered style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows (2-45) - team record repeated for every player row
$ws.Range("AD2:AD45").Value = 86
$ws.Range("AE2:AE45").Value = 76
$ws.Range("AF2:AF45").Value = 0
